function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws 'D2' '308.19'
Set-TextCell $ws 'E2' '0.56%'
Set-TextCell $ws 'D3' '40.78'
Set-TextCell $ws 'E3' '2.37%'
Set-TextCell $ws 'E4' '-0.09%'
Set-TextCell $ws 'D5' '0.07615'
Set-TextCell $ws 'E5' '-1.24%'
Set-TextCell $ws 'D6' '1.619'
Set-TextCell $ws 'E6' '-0.31%'
Set-TextCell $ws 'B7' 'BTSEToken'
Set-TextCell $ws 'C7' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell $ws 'D7' '2.448'
Set-TextCell $ws 'E7' '0.49%'
Set-TextCell $ws 'B8' 'MXToken'
Set-TextCell $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws 'D8' '0.9014'
Set-TextCell $ws 'E8' '2.34%'
Set-TextCell $ws 'D9' '0.1099'
Set-TextCell $ws 'E9' '9.68%'
Set-TextCell $ws 'D10' '0.1770'
Set-TextCell $ws 'E10' '1.41%'
Set-TextCell $ws 'D11' '0.09158'
Set-TextCell $ws 'E11' '2.45%'
Set-TextCell $ws 'D12' '0.04168'
Set-TextCell $ws 'E12' '-5.55%'
Set-TextCell $ws 'E13' '-0.56%'
Set-TextCell $ws 'D14' '0.001249'
Set-TextCell $ws 'E14' '-0.79%'
Set-TextCell $ws 'D15' '0.005882'
Set-TextCell $ws 'E15' '-0.19%'
Set-TextCell $ws 'E16' '-0.01%'
Set-TextCell $ws 'D18' '0.3292'
Set-TextCell $ws 'E18' '-0.90%'
Set-TextCell $ws 'D19' '6.603'
Set-TextCell $ws 'E19' '-5.74%'
Set-TextCell $ws 'D20' '0.1364'
Set-TextCell $ws 'E20' '1.89%'
Set-TextCell $ws 'D22' '0.04054'
Set-TextCell $ws 'E22' '-2.62%'
Set-TextCell $ws 'D23' '0.001222'
Set-TextCell $ws 'E23' '1.86%'
Set-TextCell $ws 'D24' '0.004093'
Set-TextCell $ws 'E24' '0.26%'
Set-TextCell $ws 'D25' '0.0001302'
Set-TextCell $ws 'E25' '6.59%'
Set-TextCell $ws 'D38' '0.02373'
Set-TextCell $ws 'E38' '1.54%'
Set-TextCell $ws 'D39' '0.05182'
Set-TextCell $ws 'E39' '0.89%'
Set-TextCell $ws 'D40' '0.007792'
Set-TextCell $ws 'E40' '-1.73%'
Set-TextCell $ws 'D41' '0.1299'
Set-TextCell $ws 'E41' '-1.90%'
Set-TextCell $ws 'D42' '0.006739'
Set-TextCell $ws 'E42' '6.36%'
Set-TextCell $ws 'D43' '0.001952'
Set-TextCell $ws 'E43' '0.81%'
Set-TextCell $ws 'D44' '0.007943'
Set-TextCell $ws 'E44' '-7.49%'
Set-TextCell $ws 'D45' '0.3076'
Set-TextCell $ws 'E45' '0.80%'
Set-TextCell $ws 'D46' '0.00006940'
Set-TextCell $ws 'E46' '6.55%'
Set-TextCell $ws 'E47' '-0.01%'
Set-TextCell $ws 'D48' '0.03134'
Set-TextCell $ws 'E48' '582.63%'
Set-TextCell $ws 'E49' '-39.97%'
Set-TextCell $ws 'D50' '0.00002102'
Set-TextCell $ws 'E50' '-0.01%'
Set-TextCell $ws 'E51' '-0.01%'
